$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = "'2024-03-13"
$ws.Cells.Item(2, 2).Value = "오상헬스케어"
$ws.Cells.Item(2, 3).Value = "코스닥"
$ws.Cells.Item(2, 4).Value = 198
$ws.Cells.Item(2, 5).Value = "NH"
$ws.Cells.Item(2, 6).Value = 198
$ws.Cells.Item(2, 7).Value = "-"
$ws.Cells.Item(2, 8).Value = "-"
$ws.Cells.Item(2, 9).Value = "-"
$ws.Cells.Item(2, 10).Value = "-"
$ws.Cells.Item(2, 11).Value = "대표"
$ws.Cells.Item(2, 12).Value = "-"
$ws.Cells.Item(2, 13).Value = 20000
$ws.Cells.Item(2, 14).Value = 100
$ws.Cells.Item(2, 15).Value = "2024-03-04"
$ws.Cells.Item(2, 16).Value = "2024-03-07"
$ws.Cells.Item(2, 17).Value = 742500
